# Scheduled runner update: refresh computed market/profit figures on several
# Leve-profit sheets. Values below reproduce the upstream data refresh.

$wb = $excel.ActiveWorkbook

function Set-RowValues($SheetName, $Row, $H, $I, $J, $K, $L, $M, $N) {
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range("H$Row").Value = $H
    $ws.Range("I$Row").Value = $I
    $ws.Range("J$Row").Value = $J
    $ws.Range("K$Row").Value = $K
    $ws.Range("L$Row").Value = $L
    $ws.Range("M$Row").Value = $M
    if ($null -eq $N) {
        $ws.Range("N$Row").ClearContents()
    } else {
        $ws.Range("N$Row").Value = $N
    }
}

# ALC
Set-RowValues "ALC" 28  382.3889  233.46153 769.6      233.46153 769.6      251.53847  -1739.6
Set-RowValues "ALC" 116 349175.72 4716.8335 439034.56  4716.8335 439034.56  -1274.8335 -445918.56
Set-RowValues "ALC" 125 490.2857  476.4     525        4287.599999999999 4725 -1827.599999999999 -9645
Set-RowValues "ALC" 129 889.4464  264       1009.21277 792       3027.63831 4208       -13027.63831
Set-RowValues "ALC" 137 22488314  497818.66 47620310   1493455.98 142860930 -1490905.98 -142866030
Set-RowValues "ALC" 138 1585.7931 1033.9375 2265       3101.8125 6795       2038.1875  -17075

# ARM
Set-RowValues "ARM" 102 250001060 500000600 1525       500000600 1525       -499998978 -4769
Set-RowValues "ARM" 132 3283979.8 4013447.2 1376.6     12040341.6 4129.799999999999 -12037811.6 -9189.799999999999

# BSM
Set-RowValues "BSM" 134 16273147  20217830  1329.75    60653490  3989.25    -60650955  -9059.25

# CRP
Set-RowValues "CRP" 31  18523906  22599902  25156      22599902  25156      -22599607  -25746
Set-RowValues "CRP" 34  18523906  22599902  25156      22599902  25156      -22599700  -25560
Set-RowValues "CRP" 132 7250127.5 10102643  9126.385   30307929  27379.155  -30305399  -32439.155
Set-RowValues "CRP" 134 43751536  73530750  4809483    220592250 14428449   -220589715 -14433519

# CUL
Set-RowValues "CUL" 12  64.91304  185.5     39.526318  556.5     118.578954 -383.5     -464.578954
Set-RowValues "CUL" 34  550.3823  106       901.2105   318       2703.6315  -234       -2871.6315
Set-RowValues "CUL" 136 1325      1325      0          3975      0          1125       $null

# GSM
Set-RowValues "GSM" 126 1887.4286 1129.2727 2721.4     3387.8181 8164.200000000001 -917.8181 -13104.2
Set-RowValues "GSM" 132 24549752  30596766  5983.4116  91790298  17950.2348 -91787768  -23010.2348

# LTW
Set-RowValues "LTW" 132 3390896.5 5000564   2122.8948  15001692  6368.6844  -14999162  -11428.6844

# WVR
Set-RowValues "WVR" 81  6505.9    17046.834 1988.3572  34093.668 3976.7144  -33032.668 -6098.7144
Set-RowValues "WVR" 84  6505.9    17046.834 1988.3572  170468.34 19883.572  -165164.34 -30491.572
Set-RowValues "WVR" 132 22261208  29269348  14279716   87808044  42839148   -87805514  -42844208
Set-RowValues "WVR" 136 22375284  21205038  26317166   63615114  78951498   -63612564  -78956598
